# Apply the price/volume refresh from the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/D/E hold inline text (coin name / link / price / 1h volume).
# D-column "prices" are display strings (e.g. "606.20", "63.505.08") that
# must round-trip verbatim, so they're written with a leading apostrophe
# (Excel's "force text" quote-prefix) and the resulting cell style is reset
# back to Normal so no stray number-format is left on the cell.

$ws.Range('D2').Value2 = "'" + '63.505.08'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value2 = '  -4.37%  '

$ws.Range('D3').Value2 = "'" + '3.086.51'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value2 = '  -5.44%  '

$ws.Range('D5').Value2 = "'" + '606.20'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value2 = '  -1.31%  '

$ws.Range('D6').Value2 = "'" + '144.20'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value2 = '  -8.64%  '

$ws.Range('E7').Value2 = '  +0.10%  '

$ws.Range('D8').Value2 = "'" + '3.081.98'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value2 = '  -5.55%  '

$ws.Range('E9').Value2 = '  -5.03%  '

$ws.Range('E10').Value2 = '  -8.17%  '

$ws.Range('E11').Value2 = '  -10.64%  '

$ws.Range('D12').Value2 = "'" + '0.466'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value2 = '  -6.12%  '

$ws.Range('E13').Value2 = '  -9.12%  '

$ws.Range('D14').Value2 = "'" + '34.94'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value2 = '  -10.68%  '

$ws.Range('D15').Value2 = "'" + '3.607.15'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value2 = '  -5.04%  '

$ws.Range('E16').Value2 = '  +0.98%  '

$ws.Range('D17').Value2 = "'" + '63.582.33'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value2 = '  -4.31%  '

$ws.Range('D18').Value2 = "'" + '3.092.14'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value2 = '  -5.21%  '

$ws.Range('E19').Value2 = '  -9.09%  '

$ws.Range('D20').Value2 = "'" + '471.23'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value2 = '  -6.86%  '

$ws.Range('E21').Value2 = '  -6.29%  '

$ws.Range('D22').Value2 = "'" + '0.697'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value2 = '  -7.73%  '

$ws.Range('E23').Value2 = '  -6.04%  '

$ws.Range('D24').Value2 = "'" + '13.44'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value2 = '  -8.27%  '

$ws.Range('D25').Value2 = "'" + '83.13'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value2 = '  -4.71%  '

$ws.Range('D27').Value2 = "'" + '2.76'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value2 = '  -9.21%  '

$ws.Range('E28').Value2 = '  -9.95%  '

$ws.Range('E29').Value2 = '  -11.02%  '

$ws.Range('D30').Value2 = "'" + '6.66'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value2 = '  -6.47%  '

$ws.Range('E31').Value2 = '  +0.11%  '

$ws.Range('D32').Value2 = "'" + '2.72'
$ws.Range('D32').Style = "Normal"

$ws.Range('E33').Value2 = '  -17.77%  '

$ws.Range('D34').Value2 = "'" + '25.94'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value2 = '  -7.40%  '

$ws.Range('E35').Value2 = '  -5.51%  '

$ws.Range('E36').Value2 = '  -9.47%  '

$ws.Range('D37').Value2 = "'" + '52.01'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value2 = '  -6.90%  '

$ws.Range('B38').Value2 = 'PEPE'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value2 = "'" + '0.0₃0721'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value2 = '  -8.66%  '

$ws.Range('B39').Value2 = 'Bittensor'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value2 = "'" + '456.40'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value2 = '  -8.23%  '

$ws.Range('D40').Value2 = "'" + '2.88'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value2 = '  -14.72%  '

$ws.Range('D41').Value2 = "'" + '0.0390'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value2 = '  -7.88%  '

$ws.Range('E42').Value2 = '  -8.55%  '

$ws.Range('D43').Value2 = "'" + '8.28'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value2 = '  -6.44%  '

$ws.Range('D44').Value2 = "'" + '2.812.06'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value2 = '  -6.55%  '

$ws.Range('E45').Value2 = '  -10.70%  '

$ws.Range('E46').Value2 = '  -12.02%  '

$ws.Range('B47').Value2 = 'ThetaToken'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value2 = "'" + '2.37'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value2 = '  -5.52%  '

$ws.Range('B48').Value2 = 'USDe'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value2 = "'" + '0.999'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value2 = '  -0.02%  '

$ws.Range('D49').Value2 = "'" + '25.88'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value2 = '  -10.85%  '

$ws.Range('E50').Value2 = '  -5.99%  '

$ws.Range('D51').Value2 = "'" + '117.77'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value2 = '  -2.58%  '
